# Auto-generated edit script: update cryptos list values (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text storage (avoid Excel auto-coercing numeric-looking strings to numbers),
    # then clear the format so no extra style index is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

Set-TextValue "D2" "62.363.79"
Set-TextValue "E2" "  -6.82%  "
Set-TextValue "D3" "2.916.65"
Set-TextValue "E3" "  -9.26%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.32%  "
Set-TextValue "D5" "537.03"
Set-TextValue "E5" "  -10.14%  "
Set-TextValue "D6" "130.61"
Set-TextValue "E6" "  -13.90%  "
Set-TextValue "D7" "0.998"
Set-TextValue "E7" "  -0.13%  "
Set-TextValue "D8" "2.891.18"
Set-TextValue "E8" "  -9.81%  "
Set-TextValue "D9" "0.456"
Set-TextValue "E9" "  -16.09%  "
Set-TextValue "E10" "  -18.96%  "
Set-TextValue "D11" "5.73"
Set-TextValue "E11" "  -12.07%  "
Set-TextValue "E12" "  -14.13%  "
Set-TextValue "D13" "31.63"
Set-TextValue "E13" "  -19.10%  "
Set-TextValue "D14" "0.0000200"
Set-TextValue "E14" "  -18.52%  "
Set-TextValue "D15" "3.403.07"
Set-TextValue "E15" "  -8.81%  "
Set-TextValue "D16" "62.513.30"
Set-TextValue "E16" "  -6.51%  "
Set-TextValue "E17" "  -5.75%  "
Set-TextValue "D18" "2.937.22"
Set-TextValue "E18" "  -8.63%  "
Set-TextValue "D19" "464.40"
Set-TextValue "E19" "  -12.98%  "
Set-TextValue "D20" "6.12"
Set-TextValue "E20" "  -14.65%  "
Set-TextValue "D21" "12.70"
Set-TextValue "E21" "  -15.19%  "
Set-TextValue "D22" "0.626"
Set-TextValue "E22" "  -17.98%  "
Set-TextValue "D23" "6.40"
Set-TextValue "E23" "  -19.31%  "
Set-TextValue "D24" "74.13"
Set-TextValue "E24" "  -13.21%  "
Set-TextValue "E25" "  -15.62%  "
Set-TextValue "E26" "  +0.30%  "
Set-TextValue "D27" "2.55"
Set-TextValue "E27" "  -20.14%  "
Set-TextValue "B28" "ImmutableX"
Set-TextValue "C28" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D28" "1.86"
Set-TextValue "E28" "  -15.42%  "
Set-TextValue "B29" "RenderToken"
Set-TextValue "C29" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D29" "6.88"
Set-TextValue "E29" "  -16.46%  "
Set-TextValue "D30" "24.14"
Set-TextValue "E30" "  -17.35%  "
Set-TextValue "D31" "2.34"
Set-TextValue "E31" "  -11.75%  "
Set-TextValue "B32" "FirstDigitalUSD"
Set-TextValue "C32" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  -0.03%  "
Set-TextValue "B33" "Mantle"
Set-TextValue "C33" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D33" "1.04"
Set-TextValue "E33" "  -9.51%  "
Set-TextValue "D34" "468.83"
Set-TextValue "E34" "  -14.85%  "
Set-TextValue "D35" "50.32"
Set-TextValue "E35" "  -5.95%  "
Set-TextValue "D36" "5.38"
Set-TextValue "E36" "  -18.02%  "
Set-TextValue "D37" "4.70"
Set-TextValue "E37" "  -17.81%  "
Set-TextValue "D38" "0.0384"
Set-TextValue "E38" "  -10.85%  "
Set-TextValue "D39" "0.0736"
Set-TextValue "E39" "  -15.35%  "
Set-TextValue "D40" "0.113"
Set-TextValue "E40" "  -10.12%  "
Set-TextValue "D41" "7.68"
Set-TextValue "E41" "  -18.31%  "
Set-TextValue "D42" "2.650.24"
Set-TextValue "E42" "  -9.16%  "
Set-TextValue "D43" "0.998"
Set-TextValue "E43" "  -0.24%  "
Set-TextValue "D44" "2.21"
Set-TextValue "E44" "  -17.47%  "
Set-TextValue "E45" "  -17.79%  "
Set-TextValue "D46" "111.45"
Set-TextValue "E46" "  -8.48%  "
Set-TextValue "D47" "0.0997"
Set-TextValue "E47" "  -12.70%  "
Set-TextValue "B48" "PEPE"
Set-TextValue "C48" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D48" "0.0₃0466"
Set-TextValue "E48" "  -20.65%  "
Set-TextValue "B49" "Fetch.AI"
Set-TextValue "C49" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D49" "1.76"
Set-TextValue "E49" "  -17.34%  "
Set-TextValue "D50" "21.43"
Set-TextValue "E50" "  -20.10%  "
Set-TextValue "E51" "  -5.82%  "
